$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algorithm")

# Row 7
$ws.Range("G7").Value = 1.67
$ws.Range("H7").Value = 0.033333333333333333
$ws.Range("H7").NumberFormat = "h:mm"

# Row 8
$ws.Range("G8").Value = 6.3
$ws.Range("H8").Value = 0.089583333333333334
$ws.Range("H8").NumberFormat = "h:mm"

# Row 9
$ws.Range("G9").Value = 1.45
$ws.Range("H9").Value = 0.024999999999999998
$ws.Range("H9").NumberFormat = "h:mm"

# Row 10
$ws.Range("G10").Value = 1.45
$ws.Range("H10").Value = 0.021527777777777781
$ws.Range("H10").NumberFormat = "h:mm"

# Row 11
$ws.Range("G11").Value = 2.44
$ws.Range("H11").Value = 0.037499999999999999
$ws.Range("H11").NumberFormat = "h:mm"

# Row 12
$ws.Range("G12").Value = 16.5
$ws.Range("H12").Value = 1.4756944444444444
$ws.Range("H12").NumberFormat = "[h]:mm:ss"

# Update the active selection to G9
$ws.Range("G9").Select()
